$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 367.3
$ws.Range("I53").Value = 290.33334
$ws.Range("J53").Value = 400.2857
$ws.Range("K53").Value = 290.33334
$ws.Range("L53").Value = 400.2857
$ws.Range("M53").Value = 346.66666
$ws.Range("N53").Value = -1674.2857

$ws.Range("H76").Value = 3185.7144
$ws.Range("I76").Value = 3180
$ws.Range("K76").Value = 3180
$ws.Range("M76").Value = -2865

$ws.Range("H79").Value = 3185.7144
$ws.Range("I79").Value = 3180
$ws.Range("K79").Value = 3180
$ws.Range("M79").Value = -2088

$ws.Range("H112").Value = 5276.41
$ws.Range("J112").Value = 5832
$ws.Range("L112").Value = 17496
$ws.Range("N112").Value = -19712

$ws.Range("H140").Value = 72252.37
$ws.Range("J140").Value = 72252.37
$ws.Range("L140").Value = 72252.37
$ws.Range("N140").Value = -82612.37

$ws.Range("H141").Value = 3610.7917
$ws.Range("I141").Value = 1798.375
$ws.Range("J141").Value = 7235.625
$ws.Range("K141").Value = 5395.125
$ws.Range("L141").Value = 21706.875
$ws.Range("M141").Value = -215.125
$ws.Range("N141").Value = -32066.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7156776
$ws.Range("I32").Value = 8209160
$ws.Range("J32").Value = 23950.334
$ws.Range("K32").Value = 8209160
$ws.Range("L32").Value = 23950.334
$ws.Range("M32").Value = -8208873
$ws.Range("N32").Value = -24524.334

$ws.Range("H61").Value = 6175125
$ws.Range("I61").Value = 8548585
$ws.Range("J61").Value = 4129.2666
$ws.Range("K61").Value = 8548585
$ws.Range("L61").Value = 4129.2666
$ws.Range("M61").Value = -8548373
$ws.Range("N61").Value = -4553.2666

$ws.Range("H93").Value = 53612
$ws.Range("J93").Value = 53612
$ws.Range("L93").Value = 53612
$ws.Range("N93").Value = -58604

$ws.Range("H96").Value = 82500
$ws.Range("J96").Value = 82500
$ws.Range("L96").Value = 82500
$ws.Range("N96").Value = -87992

$ws.Range("H106").Value = 54340.6
$ws.Range("J106").Value = 54340.6
$ws.Range("L106").Value = 54340.6
$ws.Range("N106").Value = -56864.6

$ws.Range("H136").Value = 6175125
$ws.Range("I136").Value = 8548585
$ws.Range("J136").Value = 4129.2666
$ws.Range("K136").Value = 25645755
$ws.Range("L136").Value = 12387.7998
$ws.Range("M136").Value = -25643205
$ws.Range("N136").Value = -17487.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 27028530
$ws.Range("I20").Value = 1436.7142
$ws.Range("J20").Value = 62501588
$ws.Range("K20").Value = 1436.7142
$ws.Range("L20").Value = 62501588
$ws.Range("M20").Value = -1189.7142
$ws.Range("N20").Value = -62502082

$ws.Range("H27").Value = 58326.668
$ws.Range("J27").Value = 58326.668
$ws.Range("L27").Value = 58326.668
$ws.Range("N27").Value = -58710.668

$ws.Range("H36").Value = 29159.857
$ws.Range("I36").Value = 509.25
$ws.Range("K36").Value = 509.25
$ws.Range("M36").Value = 24.75

$ws.Range("H69").Value = 32299.5
$ws.Range("J69").Value = 32299.5
$ws.Range("L69").Value = 32299.5
$ws.Range("N69").Value = -33921.5

$ws.Range("H72").Value = 32299.5
$ws.Range("J72").Value = 32299.5
$ws.Range("L72").Value = 96898.5
$ws.Range("N72").Value = -105010.5

$ws.Range("H75").Value = 24512.21
$ws.Range("J75").Value = 32023.715
$ws.Range("L75").Value = 32023.715
$ws.Range("N75").Value = -33895.715

$ws.Range("H78").Value = 24512.21
$ws.Range("J78").Value = 32023.715
$ws.Range("L78").Value = 96071.145
$ws.Range("N78").Value = -105431.145

$ws.Range("H86").Value = 2518.6
$ws.Range("I86").Value = 1928.6666
$ws.Range("J86").Value = 3403.5
$ws.Range("K86").Value = 1928.6666
$ws.Range("L86").Value = 3403.5
$ws.Range("M86").Value = -805.6666
$ws.Range("N86").Value = -5649.5

$ws.Range("H89").Value = 2518.6
$ws.Range("I89").Value = 1928.6666
$ws.Range("J89").Value = 3403.5
$ws.Range("K89").Value = 9643.333000000001
$ws.Range("L89").Value = 17017.5
$ws.Range("M89").Value = -4027.333000000001
$ws.Range("N89").Value = -28249.5

$ws.Range("H106").Value = 63667.75
$ws.Range("J106").Value = 63667.75
$ws.Range("L106").Value = 63667.75
$ws.Range("N106").Value = -66191.75

$ws.Range("H126").Value = 64390
$ws.Range("J126").Value = 64390
$ws.Range("L126").Value = 64390
$ws.Range("N126").Value = -74270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1778.8125
$ws.Range("I16").Value = 1764
$ws.Range("K16").Value = 1764
$ws.Range("M16").Value = -1477

$ws.Range("H28").Value = 69547.664
$ws.Range("J28").Value = 69547.664
$ws.Range("L28").Value = 69547.664
$ws.Range("N28").Value = -70037.664

$ws.Range("H31").Value = 4961.0127
$ws.Range("I31").Value = 1401.6765
$ws.Range("J31").Value = 7591.826
$ws.Range("K31").Value = 1401.6765
$ws.Range("L31").Value = 7591.826
$ws.Range("M31").Value = -1106.6765
$ws.Range("N31").Value = -8181.826

$ws.Range("H34").Value = 4961.0127
$ws.Range("I34").Value = 1401.6765
$ws.Range("J34").Value = 7591.826
$ws.Range("K34").Value = 1401.6765
$ws.Range("L34").Value = 7591.826
$ws.Range("M34").Value = -1199.6765
$ws.Range("N34").Value = -7995.826

$ws.Range("H99").Value = 2487.8333
$ws.Range("I99").Value = 2369.9
$ws.Range("J99").Value = 2533.1924
$ws.Range("K99").Value = 2369.9
$ws.Range("L99").Value = 2533.1924
$ws.Range("M99").Value = -871.9000000000001
$ws.Range("N99").Value = -5529.1924

$ws.Range("H113").Value = 1778.8125
$ws.Range("I113").Value = 1764
$ws.Range("K113").Value = 1764
$ws.Range("M113").Value = 406

$ws.Range("H126").Value = 2487.8333
$ws.Range("I126").Value = 2369.9
$ws.Range("J126").Value = 2533.1924
$ws.Range("K126").Value = 7109.700000000001
$ws.Range("L126").Value = 7599.5772
$ws.Range("M126").Value = -4639.700000000001
$ws.Range("N126").Value = -12539.5772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 139.90909
$ws.Range("J2").Value = 177.5
$ws.Range("L2").Value = 1065
$ws.Range("N2").Value = -1291

$ws.Range("H5").Value = 895.2075
$ws.Range("I5").Value = 742.44116
$ws.Range("J5").Value = 1168.579
$ws.Range("K5").Value = 2227.32348
$ws.Range("L5").Value = 3505.737
$ws.Range("M5").Value = -2115.32348
$ws.Range("N5").Value = -3729.737

$ws.Range("H131").Value = 3995.3333
$ws.Range("I131").Value = 620
$ws.Range("J131").Value = 4866.387
$ws.Range("K131").Value = 1860
$ws.Range("L131").Value = 14599.161
$ws.Range("M131").Value = 3180
$ws.Range("N131").Value = -24679.161

$ws.Range("H135").Value = 895.2075
$ws.Range("I135").Value = 742.44116
$ws.Range("J135").Value = 1168.579
$ws.Range("K135").Value = 6681.97044
$ws.Range("L135").Value = 10517.211
$ws.Range("M135").Value = -4146.97044
$ws.Range("N135").Value = -15587.211

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5060.3076
$ws.Range("I7").Value = 5338.4
$ws.Range("J7").Value = 4133.3335
$ws.Range("K7").Value = 5338.4
$ws.Range("L7").Value = 4133.3335
$ws.Range("M7").Value = -5226.4
$ws.Range("N7").Value = -4357.3335

$ws.Range("H22").Value = 15086.143
$ws.Range("I22").Value = 866.6667
$ws.Range("J22").Value = 25750.75
$ws.Range("K22").Value = 866.6667
$ws.Range("L22").Value = 25750.75
$ws.Range("M22").Value = -571.6667
$ws.Range("N22").Value = -26340.75

$ws.Range("H27").Value = 15086.143
$ws.Range("I27").Value = 866.6667
$ws.Range("J27").Value = 25750.75
$ws.Range("K27").Value = 866.6667
$ws.Range("L27").Value = 25750.75
$ws.Range("M27").Value = -759.6667
$ws.Range("N27").Value = -25964.75

$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1376

$ws.Range("H68").Value = 1875.375

$ws.Range("H71").Value = 1875.375

$ws.Range("H82").Value = 1805.8334
$ws.Range("I82").Value = 1570.2
$ws.Range("J82").Value = 2100.375
$ws.Range("K82").Value = 1570.2
$ws.Range("L82").Value = 2100.375
$ws.Range("M82").Value = -1209.2
$ws.Range("N82").Value = -2822.375

$ws.Range("H85").Value = 1805.8334
$ws.Range("I85").Value = 1570.2
$ws.Range("J85").Value = 2100.375
$ws.Range("K85").Value = 1570.2
$ws.Range("L85").Value = 2100.375
$ws.Range("M85").Value = -322.2
$ws.Range("N85").Value = -4596.375

$ws.Range("H93").Value = 12190.8
$ws.Range("I93").Value = 26450
$ws.Range("J93").Value = 2684.6667
$ws.Range("K93").Value = 26450
$ws.Range("L93").Value = 2684.6667
$ws.Range("M93").Value = -25202
$ws.Range("N93").Value = -5180.6667

$ws.Range("H126").Value = 5060.3076
$ws.Range("I126").Value = 5338.4
$ws.Range("J126").Value = 4133.3335
$ws.Range("K126").Value = 16015.2
$ws.Range("L126").Value = 12400.0005
$ws.Range("M126").Value = -13545.2
$ws.Range("N126").Value = -17340.0005

